$d = $word.ActiveDocument

$pairs = @(
    @("91×98=8918", "95×57=5415"),
    @("75×72=5400", "61×23=1403"),
    @("74×52=3848", "35×81=2835"),
    @("51×75=3825", "21×25=525"),
    @("85×75=6375", "81×31=2511"),
    @("70×28=1960", "62×40=2480"),
    @("72×61=4392", "31×12=372"),
    @("42×32=1344", "34×34=1156"),
    @("54×52=2808", "20×63=1260"),
    @("29×21=609", "60×17=1020"),
    @("16×34=544", "76×58=4408"),
    @("70×96=6720", "66×67=4422"),
    @("92×86=7912", "62×84=5208"),
    @("33×18=594", "45×74=3330"),
    @("20×84=1680", "36×95=3420"),
    @("84×46=3864", "49×96=4704"),
    @("85×50=4250", "44×52=2288"),
    @("51×18=918", "56×17=952"),
    @("95×64=6080", "52×76=3952"),
    @("79×53=4187", "49×60=2940"),
    @("32×99=3168", "89×42=3738"),
    @("48×12=576", "89×11=979"),
    @("11×11=121", "84×19=1596"),
    @("12×36=432", "96×90=8640"),
    @("84×81=6804", "39×18=702")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
